$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Kayitlar")

# A new "Kadastro App" record is appended as row 2 of the "Kayitlar"
# (Records) sheet. Every field in this workbook is stored as text (even
# the numeric-looking ones such as the record number or counts), so each
# cell is written with a literal-string formula ( ="value" ). This keeps
# the cell's stored type as a text result instead of letting Excel infer
# a number/date, while avoiding any extra number-format/style being
# introduced on the cells.
$ws.Range("A2").Formula = '="1"'
$ws.Range("B2").Formula = '="2025-07-16"'
$ws.Range("C2").Formula = '="Merkez"'
$ws.Range("D2").Formula = '="1"'
$ws.Range("E2").Formula = '="2"'
$ws.Range("F2").Formula = '="Cins D."'
$ws.Range("G2").Formula = '="Gökhan ELGÜL"'
